# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# handback step has completed for both locales (zh-cn, de-de):
#   * The "Status" column text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (shared by every cell that shows it).
#   * The "Latest Target File" (I) and "Latest Handback File" (J) columns on
#     the per-locale sheets get populated with the source doc / handback
#     xliff file names, with I turned into a hyperlink back to the source doc
#     (matching the look of the existing hyperlinks in column A).
#   * The "Latest Handback DateTime" column (K) is refreshed with the time the
#     handback report was generated.
#   * A handful of columns are widened so the newly-populated values are not
#     truncated.

$wb = $excel.ActiveWorkbook

# Helper: Excel's ColumnWidth setter only lands on a 1/6-character pixel grid,
# so round-trip through that grid to land as close as possible to the desired
# stored "width" attribute.
function Set-ColWidth($ws, $colIndex, $targetWidth) {
    $k = [Math]::Round($targetWidth * 6 - 5)
    if ($k -lt 0) { $k = 0 }
    $widthInput = $k / 6.0
    $ws.Columns.Item($colIndex).ColumnWidth = $widthInput
}

# ---------------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns (E, F)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
Set-ColWidth $overview 5 29.9777047293527
Set-ColWidth $overview 6 29.9777047293527

# Every cell showing the handoff/handback status (Overview E/F, and the
# "Status" column (C) on each per-locale sheet) needs to move in lockstep.
$statusText = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# ---------------------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de): populate handback columns
# ---------------------------------------------------------------------------
$docs = @(
    @{ Row = 2; Name = "3bc9442b-1dae-4e13-b778-38cf7e134425"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19b692ac8810dc3fed990872e8e50d469dafaa2b/e2e/3bc9442b-1dae-4e13-b778-38cf7e134425.md" },
    @{ Row = 3; Name = "56f5903a-27d0-450b-b254-6d866f433341"; Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19b692ac8810dc3fed990872e8e50d469dafaa2b/e2e/56f5903a-27d0-450b-b254-6d866f433341.md" }
)

$locales = @(
    @{ Sheet = "zh-cn"; Hash = "44e7a2c8e6e0119990240e36f52605f66f789908"; Hash2 = "e8c991cc63ce040d79262ad502beab0c730d997b"; HandbackTime = "2016-08-31 00:51:14" },
    @{ Sheet = "de-de"; Hash = "44e7a2c8e6e0119990240e36f52605f66f789908"; Hash2 = "e8c991cc63ce040d79262ad502beab0c730d997b"; HandbackTime = "2016-08-31 00:51:22" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Widen Status (C), Latest Target File (I) and Latest Handback File (J)
    Set-ColWidth $ws 3 29.9777047293527
    Set-ColWidth $ws 9 40
    Set-ColWidth $ws 10 40

    foreach ($doc in $docs) {
        $row = $doc.Row
        $mdName = "$($doc.Name).md"
        if ($row -eq 2) {
            $xlfName = "$($doc.Name).$($locale.Hash).$($locale.Sheet).xlf"
        } else {
            $xlfName = "$($doc.Name).$($locale.Hash2).$($locale.Sheet).xlf"
        }

        $targetCell = $ws.Cells.Item($row, 9)   # column I: Latest Target File
        $handbackCell = $ws.Cells.Item($row, 10) # column J: Latest Handback File
        $dateCell = $ws.Cells.Item($row, 11)     # column K: Latest Handback DateTime

        $targetCell.Value = $mdName
        $ws.Hyperlinks.Add($targetCell, $doc.Url, "", "", $mdName)
        $targetCell.Font.Underline = 2
        $targetCell.Font.Color = 15570276

        $handbackCell.Value = $xlfName

        $dateCell.Value = $locale.HandbackTime
    }
}
